$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.696.81'
$ws.Range('E2').Value = '  -0.61%  '
$ws.Range('D3').Value = '1.900.30'
$ws.Range('E3').Value = '  -0.25%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9998'
$ws.Range('E4').Value = '  -0.90%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '311.94'
$ws.Range('E5').Value = '  -1.61%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  -0.79%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4980'
$ws.Range('E7').Value = '  +2.93%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3759'
$ws.Range('E8').Value = '  -0.96%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07239'
$ws.Range('E9').Value = '  -1.79%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.95'
$ws.Range('E10').Value = '  +1.05%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.8878'
$ws.Range('E11').Value = '  -4.65%  '
$ws.Range('D12').Value = '1.956.90'
$ws.Range('E12').Value = '  +2.49%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07610'
$ws.Range('E13').Value = '  -1.80%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.437'
$ws.Range('E14').Value = '  -0.86%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '91.52'
$ws.Range('E15').Value = '  -0.40%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.003'
$ws.Range('E16').Value = '  -0.60%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008677'
$ws.Range('E17').Value = '  -1.99%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.001'
$ws.Range('E18').Value = '  -0.48%  '
$ws.Range('D19').Value = '27.702.32'
$ws.Range('E19').Value = '  -0.81%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.44'
$ws.Range('E20').Value = '  -1.54%  '
$ws.Range('E21').Value = '  -0.76%  '
$ws.Range('D22').Value = '2.216.82'
$ws.Range('E22').Value = '  +1.55%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.79'
$ws.Range('E23').Value = '  -0.90%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.567'
$ws.Range('E24').Value = '  -0.92%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '153.34'
$ws.Range('E25').Value = '  -1.57%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.844'
$ws.Range('E26').Value = '  -3.40%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.195'
$ws.Range('E27').Value = '  +3.48%  '
$ws.Range('E28').Value = '  -1.18%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '114.49'
$ws.Range('E29').Value = '  -2.35%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.820'
$ws.Range('E30').Value = '  -2.80%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08896'
$ws.Range('E31').Value = '  -0.76%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.202'
$ws.Range('E32').Value = '  -0.81%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.775'
$ws.Range('E33').Value = '  +2.69%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.221'
$ws.Range('E34').Value = '  -2.49%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7751'
$ws.Range('E35').Value = '  +1.19%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.616'
$ws.Range('E36').Value = '  +3.31%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02069'
$ws.Range('E37').Value = '  +1.08%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.048'
$ws.Range('E38').Value = '  +1.47%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.089'
$ws.Range('E39').Value = '  -0.77%  '
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5486'
$ws.Range('E40').Value = '  +0.25%  '
$ws.Range('B41').Value = 'Hedera'
$ws.Range('C41').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.05274'
$ws.Range('E41').Value = '  -0.10%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.726'
$ws.Range('E42').Value = '  -3.32%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '113.13'
$ws.Range('E43').Value = '  +3.71%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.426'
$ws.Range('E44').Value = '  -0.20%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.1505'
$ws.Range('E45').Value = '  -1.46%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4757'
$ws.Range('E46').Value = '  -1.15%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.40'
$ws.Range('E47').Value = '  -2.76%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.001'
$ws.Range('E48').Value = '  -0.79%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.609'
$ws.Range('E49').Value = '  -2.53%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '66.55'
$ws.Range('E50').Value = '  -2.07%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06005'
$ws.Range('E51').Value = '  -1.41%  '
